$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 62.583332
$ws.Range("I9").Value = 47
$ws.Range("J9").Value = 84.40000000000001
$ws.Range("K9").Value = 47
$ws.Range("L9").Value = 84.40000000000001
$ws.Range("M9").Value = 122
$ws.Range("N9").Value = -422.4

$ws.Range("H40").Value = 1181.9778
$ws.Range("I40").Value = 1115.9062
$ws.Range("J40").Value = 1344.6154
$ws.Range("K40").Value = 1115.9062
$ws.Range("L40").Value = 1344.6154
$ws.Range("M40").Value = -940.9061999999999
$ws.Range("N40").Value = -1694.6154

$ws.Range("H51").Value = 3232.5557
$ws.Range("I51").Value = 2515.1667
$ws.Range("J51").Value = 4667.3335
$ws.Range("K51").Value = 2515.1667
$ws.Range("L51").Value = 4667.3335
$ws.Range("M51").Value = -2031.1667
$ws.Range("N51").Value = -5635.3335

$ws.Range("H86").Value = 2962.1904
$ws.Range("I86").Value = 3121.7856
$ws.Range("J86").Value = 2882.3928
$ws.Range("K86").Value = 3121.7856
$ws.Range("L86").Value = 2882.3928
$ws.Range("M86").Value = -1998.7856
$ws.Range("N86").Value = -5128.3928

$ws.Range("H87").Value = 14253.571
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H89").Value = 2962.1904
$ws.Range("I89").Value = 3121.7856
$ws.Range("J89").Value = 2882.3928
$ws.Range("K89").Value = 15608.928
$ws.Range("L89").Value = 14411.964
$ws.Range("M89").Value = -9992.928
$ws.Range("N89").Value = -25643.964

$ws.Range("H90").Value = 14253.571
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H100").Value = 2095.1538
$ws.Range("I100").Value = 1567
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1567
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -1026
$ws.Range("N100").Value = -6082

$ws.Range("H107").Value = 761.6667
$ws.Range("I107").Value = 782.94116
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 782.94116
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1137.05884
$ws.Range("N107").Value = -4240

$ws.Range("H109").Value = 34547.2
$ws.Range("J109").Value = 34547.2
$ws.Range("L109").Value = 34547.2
$ws.Range("N109").Value = -37321.2

$ws.Range("H113").Value = 2381.6667
$ws.Range("I113").Value = 2119.6875
$ws.Range("J113").Value = 3220
$ws.Range("K113").Value = 2119.6875
$ws.Range("L113").Value = 3220
$ws.Range("M113").Value = 1134.3125
$ws.Range("N113").Value = -9728

$ws.Range("H129").Value = 847.8125
$ws.Range("J129").Value = 929.1429000000001
$ws.Range("L129").Value = 2787.4287
$ws.Range("N129").Value = -12787.4287

$ws.Range("H132").Value = 2977471
$ws.Range("I132").Value = 2977471
$ws.Range("K132").Value = 8932413
$ws.Range("M132").Value = -8929883

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3444.23
$ws.Range("I32").Value = 2834.4268
$ws.Range("J32").Value = 6222.222
$ws.Range("K32").Value = 2834.4268
$ws.Range("L32").Value = 6222.222
$ws.Range("M32").Value = -2547.4268
$ws.Range("N32").Value = -6796.222

$ws.Range("H61").Value = 2961.1177
$ws.Range("I61").Value = 1507.8948
$ws.Range("J61").Value = 4801.8667
$ws.Range("K61").Value = 1507.8948
$ws.Range("L61").Value = 4801.8667
$ws.Range("M61").Value = -1295.8948
$ws.Range("N61").Value = -5225.8667

$ws.Range("H102").Value = 4193.625
$ws.Range("I102").Value = 2611.762
$ws.Range("K102").Value = 2611.762
$ws.Range("M102").Value = -989.7620000000002

$ws.Range("H122").Value = 3066.375
$ws.Range("I122").Value = 2468.7144
$ws.Range("K122").Value = 7406.1432
$ws.Range("M122").Value = -4956.1432

$ws.Range("H136").Value = 2961.1177
$ws.Range("I136").Value = 1507.8948
$ws.Range("J136").Value = 4801.8667
$ws.Range("K136").Value = 4523.6844
$ws.Range("L136").Value = 14405.6001
$ws.Range("M136").Value = -1973.6844
$ws.Range("N136").Value = -19505.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2229.7827
$ws.Range("I20").Value = 2134.2222
$ws.Range("J20").Value = 2291.2144
$ws.Range("K20").Value = 2134.2222
$ws.Range("L20").Value = 2291.2144
$ws.Range("M20").Value = -1887.2222
$ws.Range("N20").Value = -2785.2144

$ws.Range("H105").Value = 1697.5714
$ws.Range("I105").Value = 1597
$ws.Range("J105").Value = 2125
$ws.Range("K105").Value = 1597
$ws.Range("L105").Value = 2125
$ws.Range("M105").Value = 150
$ws.Range("N105").Value = -5619

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3620.9473
$ws.Range("I31").Value = 3438.3684
$ws.Range("J31").Value = 3803.5264
$ws.Range("K31").Value = 3438.3684
$ws.Range("L31").Value = 3803.5264
$ws.Range("M31").Value = -3143.3684
$ws.Range("N31").Value = -4393.526400000001

$ws.Range("H34").Value = 3620.9473
$ws.Range("I34").Value = 3438.3684
$ws.Range("J34").Value = 3803.5264
$ws.Range("K34").Value = 3438.3684
$ws.Range("L34").Value = 3803.5264
$ws.Range("M34").Value = -3236.3684
$ws.Range("N34").Value = -4207.526400000001

$ws.Range("H105").Value = 1572.8182
$ws.Range("I105").Value = 1572.8182
$ws.Range("K105").Value = 1572.8182
$ws.Range("M105").Value = 174.1818000000001

$ws.Range("H138").Value = 58000
$ws.Range("J138").Value = 58000
$ws.Range("L138").Value = 58000
$ws.Range("N138").Value = -68280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1465.8948
$ws.Range("I34").Value = 591.8570999999999
$ws.Range("J34").Value = 1975.75
$ws.Range("K34").Value = 1775.5713
$ws.Range("L34").Value = 5927.25
$ws.Range("M34").Value = -1691.5713
$ws.Range("N34").Value = -6095.25

$ws.Range("H39").Value = 1357.8667
$ws.Range("J39").Value = 1459.0769
$ws.Range("L39").Value = 4377.2307
$ws.Range("N39").Value = -4965.2307

$ws.Range("H55").Value = 2891.4285
$ws.Range("I55").Value = 1500
$ws.Range("J55").Value = 2998.4614
$ws.Range("K55").Value = 4500
$ws.Range("L55").Value = 8995.3842
$ws.Range("M55").Value = -4323
$ws.Range("N55").Value = -9349.3842

$ws.Range("H109").Value = 3338.3333
$ws.Range("I109").Value = 1136.8334
$ws.Range("J109").Value = 4439.0835
$ws.Range("K109").Value = 3410.5002
$ws.Range("L109").Value = 13317.2505
$ws.Range("M109").Value = -2370.5002
$ws.Range("N109").Value = -15397.2505

$ws.Range("H113").Value = 1219.95
$ws.Range("I113").Value = 4126
$ws.Range("J113").Value = 493.4375
$ws.Range("K113").Value = 12378
$ws.Range("L113").Value = 1480.3125
$ws.Range("M113").Value = -10208
$ws.Range("N113").Value = -5820.3125

$ws.Range("H131").Value = 833.35596
$ws.Range("I131").Value = 432
$ws.Range("J131").Value = 1009.561
$ws.Range("K131").Value = 1296
$ws.Range("L131").Value = 3028.683
$ws.Range("M131").Value = 3744
$ws.Range("N131").Value = -13108.683

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1319.85
$ws.Range("I102").Value = 1223.7931
$ws.Range("J102").Value = 1573.091
$ws.Range("K102").Value = 1223.7931
$ws.Range("L102").Value = 1573.091
$ws.Range("M102").Value = 398.2068999999999
$ws.Range("N102").Value = -4817.091

$ws.Range("H132").Value = 3343.5
$ws.Range("I132").Value = 3182
$ws.Range("J132").Value = 3458.8572
$ws.Range("K132").Value = 9546
$ws.Range("L132").Value = 10376.5716
$ws.Range("M132").Value = -7016
$ws.Range("N132").Value = -15436.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2643.6365
$ws.Range("I7").Value = 1526.25
$ws.Range("J7").Value = 3282.1428
$ws.Range("K7").Value = 1526.25
$ws.Range("L7").Value = 3282.1428
$ws.Range("M7").Value = -1414.25
$ws.Range("N7").Value = -3506.1428

$ws.Range("H40").Value = 2715.65
$ws.Range("I40").Value = 2437.1538
$ws.Range("J40").Value = 3232.8572
$ws.Range("K40").Value = 2437.1538
$ws.Range("L40").Value = 3232.8572
$ws.Range("M40").Value = -2301.1538
$ws.Range("N40").Value = -3504.8572

$ws.Range("H59").Value = 25132.666
$ws.Range("J59").Value = 25132.666
$ws.Range("L59").Value = 25132.666
$ws.Range("N59").Value = -26440.666

$ws.Range("H111").Value = 33000
$ws.Range("J111").Value = 33000
$ws.Range("L111").Value = 33000
$ws.Range("N111").Value = -41180

$ws.Range("H126").Value = 2643.6365
$ws.Range("I126").Value = 1526.25
$ws.Range("J126").Value = 3282.1428
$ws.Range("K126").Value = 4578.75
$ws.Range("L126").Value = 9846.428400000001
$ws.Range("M126").Value = -2108.75
$ws.Range("N126").Value = -14786.4284

$ws.Range("H132").Value = 4061.9143
$ws.Range("I132").Value = 2700.681
$ws.Range("K132").Value = 8102.043
$ws.Range("M132").Value = -5572.043

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 19000
$ws.Range("I43").Value = 19000
$ws.Range("J43").Value = 19000
$ws.Range("K43").Value = 19000
$ws.Range("L43").Value = 19000
$ws.Range("M43").Value = -18851
$ws.Range("N43").Value = -19298

$ws.Range("H100").Value = 1252.75
$ws.Range("I100").Value = 1252.75
$ws.Range("K100").Value = 2505.5
$ws.Range("M100").Value = -1964.5

$ws.Range("H101").Value = 21680.4
$ws.Range("J101").Value = 21680.4
$ws.Range("L101").Value = 21680.4
$ws.Range("N101").Value = -28170.4

$ws.Range("H107").Value = 327.1
$ws.Range("I107").Value = 306.2857
$ws.Range("J107").Value = 375.66666
$ws.Range("K107").Value = 918.8571000000001
$ws.Range("L107").Value = 1126.99998
$ws.Range("M107").Value = 1001.1429
$ws.Range("N107").Value = -4966.999980000001

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H126").Value = 667510.6
$ws.Range("I126").Value = 714904.2
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 2144712.6
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -2142242.6
$ws.Range("N126").Value = -16940
